# NATMI Ncam1-Robo1 LR-pair table refreshed with updated TPM-derived expression stats.
# Ligand (Ncam1) columns E:J are keyed by "Sending cluster" (col A); receptor (Robo1)
# columns K:P are keyed by "Target cluster" (col D); edge columns Q:T = ligand x receptor
# (Q = G*M, R = H*N, S = I*O, T = J*P). Only the changed numeric cells are rewritten below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.208684666666667
$ws.Range("H2").Value = 3.626054
$ws.Range("I2").Value = 0.01462795763842055
$ws.Range("J2").Value = 0.01462795763842055
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.312815
$ws.Range("N2").Value = 0.938445
$ws.Range("O2").Value = 0.0082131704949067
$ws.Range("P2").Value = 0.0082131704949067
$ws.Range("Q2").Value = 0.3780946940033333
$ws.Range("R2").Value = 3.40285224603
$ws.Range("S2").Value = 0.0001201419100766208
$ws.Range("T2").Value = 0.0001201419100766208
# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.208684666666667
$ws.Range("H3").Value = 3.626054
$ws.Range("I3").Value = 0.01462795763842055
$ws.Range("J3").Value = 0.01462795763842055
$ws.Range("O3").Value = 0.6524076620340182
$ws.Range("P3").Value = 0.6524076620340182
$ws.Range("Q3").Value = 30.03369715692044
$ws.Range("R3").Value = 270.303274412284
$ws.Range("S3").Value = 0.009543391643214614
$ws.Range("T3").Value = 0.009543391643214614
# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.208684666666667
$ws.Range("H4").Value = 3.626054
$ws.Range("I4").Value = 0.01462795763842055
$ws.Range("J4").Value = 0.01462795763842055
$ws.Range("M4").Value = 12.866992
$ws.Range("N4").Value = 38.600976
$ws.Range("O4").Value = 0.3378316226926476
$ws.Range("P4").Value = 0.3378316226926476
$ws.Range("Q4").Value = 15.55213593652267
$ws.Range("R4").Value = 139.969223428704
$ws.Range("S4").Value = 0.004941786665666925
$ws.Range("T4").Value = 0.004941786665666925
# Row 5: ECs -> Resolving-Mac
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.208684666666667
$ws.Range("H5").Value = 3.626054
$ws.Range("I5").Value = 0.01462795763842055
$ws.Range("J5").Value = 0.01462795763842055
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05894133333333334
$ws.Range("N5").Value = 0.176824
$ws.Range("O5").Value = 0.001547544778427486
$ws.Range("P5").Value = 0.001547544778427486
$ws.Range("Q5").Value = 0.07124148583288889
$ws.Range("R5").Value = 0.641173372496
$ws.Range("S5").Value = 0.00002263741946239619
$ws.Range("T5").Value = 0.00002263741946239619
# Row 6: FAPs -> ECs
$ws.Range("I6").Value = 0.0626664797952065
$ws.Range("J6").Value = 0.06266647979520648
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.312815
$ws.Range("N6").Value = 0.938445
$ws.Range("O6").Value = 0.0082131704949067
$ws.Range("P6").Value = 0.0082131704949067
$ws.Range("Q6").Value = 1.619765662993334
$ws.Range("R6").Value = 14.57789096694
$ws.Range("S6").Value = 0.0005146904828736568
$ws.Range("T6").Value = 0.0005146904828736567
# Row 7: FAPs -> FAPs
$ws.Range("I7").Value = 0.0626664797952065
$ws.Range("J7").Value = 0.06266647979520648
$ws.Range("O7").Value = 0.6524076620340182
$ws.Range("P7").Value = 0.6524076620340182
$ws.Range("S7").Value = 0.04088409157109271
$ws.Range("T7").Value = 0.04088409157109271
# Row 8: FAPs -> MuSCs
$ws.Range("I8").Value = 0.0626664797952065
$ws.Range("J8").Value = 0.06266647979520648
$ws.Range("M8").Value = 12.866992
$ws.Range("N8").Value = 38.600976
$ws.Range("O8").Value = 0.3378316226926476
$ws.Range("P8").Value = 0.3378316226926476
$ws.Range("Q8").Value = 66.62567916375468
$ws.Range("R8").Value = 599.6311124737921
$ws.Range("S8").Value = 0.02117071855765063
$ws.Range("T8").Value = 0.02117071855765062
# Row 9: FAPs -> Resolving-Mac
$ws.Range("I9").Value = 0.0626664797952065
$ws.Range("J9").Value = 0.06266647979520648
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.05894133333333334
$ws.Range("N9").Value = 0.176824
$ws.Range("O9").Value = 0.001547544778427486
$ws.Range("P9").Value = 0.001547544778427486
$ws.Range("Q9").Value = 0.3052000315342223
$ws.Range("R9").Value = 2.746800283808
$ws.Range("S9").Value = 0.00009697918358950338
$ws.Range("T9").Value = 0.00009697918358950335
# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 76.16218566666667
$ws.Range("H10").Value = 228.486557
$ws.Range("I10").Value = 0.9217434921665711
$ws.Range("J10").Value = 0.921743492166571
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.312815
$ws.Range("N10").Value = 0.938445
$ws.Range("O10").Value = 0.0082131704949067
$ws.Range("P10").Value = 0.0082131704949067
$ws.Range("Q10").Value = 23.82467410931834
$ws.Range("R10").Value = 214.422066983865
$ws.Range("S10").Value = 0.007570436453734747
$ws.Range("T10").Value = 0.007570436453734746
# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 76.16218566666667
$ws.Range("H11").Value = 228.486557
$ws.Range("I11").Value = 0.9217434921665711
$ws.Range("J11").Value = 0.921743492166571
$ws.Range("O11").Value = 0.6524076620340182
$ws.Range("P11").Value = 0.6524076620340182
$ws.Range("Q11").Value = 1892.496928442169
$ws.Range("R11").Value = 17032.47235597952
$ws.Range("S11").Value = 0.6013525167194641
$ws.Range("T11").Value = 0.601352516719464
# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 76.16218566666667
$ws.Range("H12").Value = 228.486557
$ws.Range("I12").Value = 0.9217434921665711
$ws.Range("J12").Value = 0.921743492166571
$ws.Range("M12").Value = 12.866992
$ws.Range("N12").Value = 38.600976
$ws.Range("O12").Value = 0.3378316226926476
$ws.Range("P12").Value = 0.3378316226926476
$ws.Range("Q12").Value = 979.9782336755148
$ws.Range("R12").Value = 8819.804103079632
$ws.Range("S12").Value = 0.3113940996650205
$ws.Range("T12").Value = 0.3113940996650204
# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 76.16218566666667
$ws.Range("H13").Value = 228.486557
$ws.Range("I13").Value = 0.9217434921665711
$ws.Range("J13").Value = 0.921743492166571
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.05894133333333334
$ws.Range("N13").Value = 0.176824
$ws.Range("O13").Value = 0.001547544778427486
$ws.Range("P13").Value = 0.001547544778427486
$ws.Range("Q13").Value = 4.489100772774223
$ws.Range("R13").Value = 40.40190695496801
$ws.Range("S13").Value = 0.001426439328351894
$ws.Range("T13").Value = 0.001426439328351893
# Row 14: Resolving-Mac -> ECs
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.07949433333333333
$ws.Range("H14").Value = 0.238483
$ws.Range("I14").Value = 0.0009620703998019471
$ws.Range("J14").Value = 0.000962070399801947
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.312815
$ws.Range("N14").Value = 0.938445
$ws.Range("O14").Value = 0.0082131704949067
$ws.Range("P14").Value = 0.0082131704949067
$ws.Range("Q14").Value = 0.02486701988166667
$ws.Range("R14").Value = 0.223803178935
$ws.Range("S14").Value = 0.000007901648221676444
$ws.Range("T14").Value = 0.000007901648221676442
# Row 15: Resolving-Mac -> FAPs
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.07949433333333333
$ws.Range("H15").Value = 0.238483
$ws.Range("I15").Value = 0.0009620703998019471
$ws.Range("J15").Value = 0.000962070399801947
$ws.Range("O15").Value = 0.6524076620340182
$ws.Range("P15").Value = 0.6524076620340182
$ws.Range("Q15").Value = 1.975294962257556
$ws.Range("R15").Value = 17.777654660318
$ws.Range("S15").Value = 0.0006276621002469215
$ws.Range("T15").Value = 0.0006276621002469214
# Row 16: Resolving-Mac -> MuSCs
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.07949433333333333
$ws.Range("H16").Value = 0.238483
$ws.Range("I16").Value = 0.0009620703998019471
$ws.Range("J16").Value = 0.000962070399801947
$ws.Range("M16").Value = 12.866992
$ws.Range("N16").Value = 38.600976
$ws.Range("O16").Value = 0.3378316226926476
$ws.Range("P16").Value = 0.3378316226926476
$ws.Range("Q16").Value = 1.022852951045333
$ws.Range("R16").Value = 9.205676559408001
$ws.Range("S16").Value = 0.000325017804309656
$ws.Range("T16").Value = 0.000325017804309656
# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.07949433333333333
$ws.Range("H17").Value = 0.238483
$ws.Range("I17").Value = 0.0009620703998019471
$ws.Range("J17").Value = 0.000962070399801947
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.05894133333333334
$ws.Range("N17").Value = 0.176824
$ws.Range("O17").Value = 0.001547544778427486
$ws.Range("P17").Value = 0.001547544778427486
$ws.Range("Q17").Value = 0.004685501999111111
$ws.Range("R17").Value = 0.042169517992
$ws.Range("S17").Value = 0.000001488847023693147
$ws.Range("T17").Value = 0.000001488847023693147
